$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh: update D (Price) and E (Volume 1h) columns.
# Values that look like plain numbers (e.g. "581.97") must be forced to remain
# TEXT cells (matching the sheet data source), since a bare .Value assignment
# would otherwise have Excel auto-coerce them into numeric cells.

$ws.Range("D2").Value = '63.832.03'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.756.61'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.97'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.73'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.65%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.97'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -12.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.392'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").Value = '3.249.10'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.98'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("D15").Value = '63.798.85'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").Value = '2.767.43'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.26'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.88%  '
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '362.03'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.88'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.566'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.992'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.60'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.172'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.69'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").Value = '0.0₃0936'
$ws.Range("E28").Value = '  +6.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.01'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.24'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.60'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.03'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +5.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.56'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.20'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.15'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +8.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '331.95'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.77'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.05'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0601'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.98'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.643'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0258'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.50'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.96%  '
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("E51").Value = '  +1.07%  '
